$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of transfer data for the 6/19 hatchery day shipment
$ws.Range("A3").Value = "Ariana Huffmyer, Steven Roberts; University of Washington"
$ws.Range("B3").Value = "Crassostrea (Magallana) gigas"
$ws.Range("D3").Value = "Juen 25 2024"
$ws.Range("E3").Value = "Goose Point Oyster Co. "
$ws.Range("F3").Value = "7081 Niawiakum St Hwy #101, Bay Center, WA 98527"
$ws.Range("G3").Value = "Scientific research, not for consumption"

# Quantity shipped is unknown/left blank, highlight it in yellow
$ws.Range("C3").Interior.Color = 65535

# Restore the view back to A1 and select F11
$ws.Range("F11").Select()
